$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# card[3] ("From sale of stock you get $50") now pays out $50 (previously blank)
$ws.Range("E5").Value = 50

# card[5] ("Go to Jail") gets an "advance" function to position 10
$ws.Range("C7").Value = "advance"
$ws.Range("D7").Value = 10

# Update the saved selection to match the authored change
$ws.Range("E12").Select()
